$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "26.547.09"
$rng.Style = "Normal"
$ws.Range("E2").Value = "  +6.68%  "
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "1.721.72"
$rng.Style = "Normal"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").Value = "  -0.26%  "
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "333.07"
$rng.Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("E6").Value = "  -0.14%  "
$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = "0.3704"
$rng.Style = "Normal"
$ws.Range("E7").Value = "  +1.68%  "
$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = "48.27"
$rng.Style = "Normal"
$ws.Range("E8").Value = "  +2.18%  "
$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = "0.3356"
$rng.Style = "Normal"
$ws.Range("E9").Value = "  +2.61%  "
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = "1.184"
$rng.Style = "Normal"
$ws.Range("E10").Value = "  +3.73%  "
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "0.07390"
$rng.Style = "Normal"
$ws.Range("E11").Value = "  +4.24%  "
$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = "0.9993"
$rng.Style = "Normal"
$ws.Range("E12").Value = "  -0.25%  "
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "6.393"
$rng.Style = "Normal"
$ws.Range("E13").Value = "  +5.14%  "
$ws.Range("E14").Value = "  +2.66%  "
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "7.051"
$rng.Style = "Normal"
$ws.Range("E15").Value = "  +6.49%  "
$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = "1.720.38"
$rng.Style = "Normal"
$ws.Range("E16").Value = "  +3.65%  "
$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = "0.00001068"
$rng.Style = "Normal"
$ws.Range("E17").Value = "  +1.89%  "
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "0.06620"
$rng.Style = "Normal"
$ws.Range("E18").Value = "  -1.12%  "
$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = "82.07"
$rng.Style = "Normal"
$ws.Range("E19").Value = "  +4.23%  "
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = "1.000"
$rng.Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "16.52"
$rng.Style = "Normal"
$ws.Range("E21").Value = "  +4.25%  "
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = "6.138"
$rng.Style = "Normal"
$ws.Range("E22").Value = "  +3.43%  "
$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = "12.78"
$rng.Style = "Normal"
$ws.Range("E23").Value = "  +1.09%  "
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = "26.486.25"
$rng.Style = "Normal"
$ws.Range("E24").Value = "  +6.56%  "
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = "2.430"
$rng.Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = "2.392"
$rng.Style = "Normal"
$ws.Range("E26").Value = "  -1.84%  "
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = "1.397"
$rng.Style = "Normal"
$ws.Range("E27").Value = "  +18.99%  "
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = "152.09"
$rng.Style = "Normal"
$ws.Range("E28").Value = "  +1.44%  "
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = "19.34"
$rng.Style = "Normal"
$ws.Range("E29").Value = "  +3.33%  "
$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = "1.911.08"
$rng.Style = "Normal"
$ws.Range("E30").Value = "  +3.65%  "
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = "130.75"
$rng.Style = "Normal"
$ws.Range("E31").Value = "  +3.80%  "
$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = "4.121"
$rng.Style = "Normal"
$ws.Range("E32").Value = "  +1.11%  "
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "5.936"
$rng.Style = "Normal"
$ws.Range("E33").Value = "  +4.26%  "
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = "0.08607"
$rng.Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = "1.702"
$rng.Style = "Normal"
$ws.Range("E35").Value = "  +2.93%  "
$rng = $ws.Range("D36")
$rng.NumberFormat = "@"
$rng.Value = "12.64"
$rng.Style = "Normal"
$ws.Range("E36").Value = "  +3.56%  "
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = "5.331"
$rng.Style = "Normal"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("E38").Value = "  +1.83%  "
$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = "0.2154"
$rng.Style = "Normal"
$ws.Range("E39").Value = "  +3.42%  "
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = "0.06176"
$rng.Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = "8.438"
$rng.Style = "Normal"
$ws.Range("E41").Value = "  +1.62%  "
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = "1.221"
$rng.Style = "Normal"
$ws.Range("E42").Value = "  -3.23%  "
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = "0.6158"
$rng.Style = "Normal"
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("E44").Value = "  -0.09%  "
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = "14.15"
$rng.Style = "Normal"
$ws.Range("E45").Value = "  +4.81%  "
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = "3.905"
$rng.Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = "0.5952"
$rng.Style = "Normal"
$ws.Range("E47").Value = "  +5.09%  "
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "127.98"
$rng.Style = "Normal"
$ws.Range("E48").Value = "  +1.61%  "
$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = "2.037"
$rng.Style = "Normal"
$ws.Range("E49").Value = "  +4.02%  "
$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = "0.07158"
$rng.Style = "Normal"
$ws.Range("E50").Value = "  +2.66%  "
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = "76.83"
$rng.Style = "Normal"
$ws.Range("E51").Value = "  +2.27%  "
